$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values (column B) for rows 2-10
$ws.Range("B2").Value = 151
$ws.Range("B3").Value = 148
$ws.Range("B4").Value = 116
$ws.Range("B5").Value = 115
$ws.Range("B6").Value = 108
$ws.Range("B7").Value = 88
$ws.Range("B8").Value = 85
$ws.Range("B9").Value = 78
$ws.Range("B10").Value = 72

# Update empadronador names (column A) for rows 7-9 to reflect reordering
$ws.Range("A7").Value = "SAUCEDO CABRERA CARLOS ALEXANDER"
$ws.Range("A8").Value = "BECERRA ASMAT CAROL STEFANY"
$ws.Range("A9").Value = "PACHECO ALISON"
